$d = $word.ActiveDocument

# --- Remove the stray _GoBack bookmark from the signature paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Locate the paragraph ending "...pozostaje anonimowy. " ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*pozostaje anonimowy.*") {
        $target = $p
    }
}

# Append the new sentence as four separate runs, matching the diff:
#   "Mozliwosc rezerwacji " | "dla niezarejestrowanych" | "_GoBack" bookmark | " " | "klientow."
$r1 = $d.Range($target.Range.End - 1, $target.Range.End - 1)
$r1.InsertAfter("Mozliwosc rezerwacji ")

$r2 = $d.Range($target.Range.End - 1, $target.Range.End - 1)
$r2.InsertAfter("dla niezarejestrowanych")

$r3 = $d.Range($target.Range.End - 1, $target.Range.End - 1)
$r3.InsertAfter(" ")

$r4 = $d.Range($target.Range.End - 1, $target.Range.End - 1)
$r4.InsertAfter("klientow.")

# Re-insert the _GoBack bookmark between "dla niezarejestrowanych" and " klientow."
# (its position is no longer the paragraph's final character, so it lands correctly)
$bmPos = $target.Range.End - 1 - ("klientow.").Length - (" ").Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
